# Sprint 2 - "add de aplicação funcional"
# Adds 9 grey "Retângulo" rectangles (functional-application mock-up
# highlight boxes) onto the last slide ("OBRIGADO!") of the deck, on top
# of the existing screenshot picture.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)

# EMU -> point conversion helper (PowerPoint COM works in points; the
# OOXML stores English Metric Units, 12700 EMU == 1 point).
function Pt([double]$emu) { return $emu / 12700.0 }

$msoShapeRectangle = 1
$ppAlignCenter = 2
$msoAnchorMiddle = 3

$rects = @(
    @{ Name = "Retângulo 1";  X = 1317171; Y = 3743325; CX = 1643743; CY = 262618; Fill = 0x7F7F7F },
    @{ Name = "Retângulo 3";  X = 3429000; Y = 3820659; CX = 98425;   CY = 138566; Fill = 0xA6A6A6 },
    @{ Name = "Retângulo 8";  X = 3492500; Y = 3646232; CX = 561975;  CY = 45719;  Fill = 0xA6A6A6 },
    @{ Name = "Retângulo 10"; X = 6924675; Y = 2572884; CX = 98425;   CY = 138566; Fill = 0xA6A6A6 },
    @{ Name = "Retângulo 11"; X = 1498600; Y = 4102498; CX = 98425;   CY = 138566; Fill = 0xA6A6A6 },
    @{ Name = "Retângulo 12"; X = 5811157; Y = 1839204; CX = 561975;  CY = 45719;  Fill = 0x7F7F7F },
    @{ Name = "Retângulo 13"; X = 7471341; Y = 1538599; CX = 561975;  CY = 45719;  Fill = 0x7F7F7F },
    @{ Name = "Retângulo 14"; X = 6692899; Y = 1538599; CX = 701676;  CY = 45719;  Fill = 0x7F7F7F },
    @{ Name = "Retângulo 15"; X = 5460318; Y = 3948390; CX = 743631;  CY = 45719;  Fill = 0x7F7F7F }
)

foreach ($r in $rects) {
    $shp = $s.Shapes.AddShape($msoShapeRectangle, (Pt $r.X), (Pt $r.Y), (Pt $r.CX), (Pt $r.CY))
    $shp.Name = $r.Name

    $shp.Fill.ForeColor.RGB = $r.Fill
    $shp.Line.ForeColor.RGB = $r.Fill

    $shp.TextFrame.VerticalAnchor = $msoAnchorMiddle
    $shp.TextFrame.TextRange.ParagraphFormat.Alignment = $ppAlignCenter
}
